$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the divider/asterisk rows' content while keeping the cell style.
# These cells previously held the "****...****" separator string.
$rows = @(38, 76, 114, 152, 191, 229)
foreach ($r in $rows) {
    $ws.Range("A$r").ClearContents()
}

# Update the sheet view to match scroll/selection state after the edit.
$excel.ActiveWindow.ScrollRow = 212
$ws.Range("A229").Select()
